{"js": "const body = context.document.body;\n\nconst replacements = [\n  {\n    oldText: \"Ativa\u00e7\u00e3o: 01/01/2023\",\n    newText: \"Ativa\u00e7\u00e3o: 01/01/2024\"\n  },\n  {\n    oldText: \"Organiza\u00e7\u00e3o e o formalismo do desenvolvimento do trabalho cient\u00edfico. T\u00e9cnicas de reda\u00e7\u00e3o cient\u00edfica, uso de ferramentas de busca, refer\u00eancias bibliogr\u00e1ficas e estruturas formais de divulga\u00e7\u00e3o cient\u00edfica. Desenvolvimento de um tema de pesquisa individual, com o formato de um trabalho de Inicia\u00e7\u00e3o Cient\u00edfica, sob a orienta\u00e7\u00e3o de um professor ou pesquisador autorizado pela Comiss\u00e3o de Curso. Entrega e apresenta\u00e7\u00e3o de monografia no final da disciplina.\",\n    newText: \"Organiza\u00e7\u00e3o e o formalismo do desenvolvimento do trabalho cient\u00edfico ou projeto de engenharia. T\u00e9cnicas de reda\u00e7\u00e3o cient\u00edfica, uso de ferramentas de busca, refer\u00eancias bibliogr\u00e1ficas e estruturas formais de divulga\u00e7\u00e3o cient\u00edfica. Desenvolvimento de um tema de pesquia ou projeto de engenharia, com o formato de um trabalho de inicia\u00e7\u00e3o cient\u00edfica, sob a orienta\u00e7\u00e3o de um professor ou pesquisador autorizado pela Comiss\u00e3o de Curso. Entrega e apresenta\u00e7\u00e3o de documento t\u00e9cnico no final da disciplina.\"\n  },\n  {\n    oldText: \"Organization and formalism of the development of scientific work. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project at the end of the course.\",\n    newText: \"Organization and formalism of the development of scientific work or engineering project. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research or engineering project topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project document at the end of the course.\"\n  },\n  {\n    oldText: \"Aulas expositivas, reuni\u00f5es com professor orientador, desenvolvimento de projeto de pesquisa e elabora\u00e7\u00e3o de projeto de pesquisa.\",\n    newText: \"Aulas expositivas, reuni\u00f5es com professor orientador, desenvolvimento de projeto de pesquisa e/ou engenharia e elabora\u00e7\u00e3o de projeto de pesquisa e/ou engenharia. Visitas t\u00e9cnicas em institutos ou empresas da \u00e1rea cient\u00edfica ou de engenharia.\"\n  }\n];\n\nfor (const rep of replacements) {\n  const found = body.search(rep.oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(rep.newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1  # wdFindContinue\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2023\" \"Ativa\u00e7\u00e3o: 01/01/2024\"\n\nReplace-Text \"Organiza\u00e7\u00e3o e o formalismo do desenvolvimento do trabalho cient\u00edfico. T\u00e9cnicas de reda\u00e7\u00e3o cient\u00edfica, uso de ferramentas de busca, refer\u00eancias bibliogr\u00e1ficas e estruturas formais de divulga\u00e7\u00e3o cient\u00edfica. Desenvolvimento de um tema de pesquisa individual, com o formato de um trabalho de Inicia\u00e7\u00e3o Cient\u00edfica, sob a orienta\u00e7\u00e3o de um professor ou pesquisador autorizado pela Comiss\u00e3o de Curso. Entrega e apresenta\u00e7\u00e3o de monografia no final da disciplina.\" \"Organiza\u00e7\u00e3o e o formalismo do desenvolvimento do trabalho cient\u00edfico ou projeto de engenharia. T\u00e9cnicas de reda\u00e7\u00e3o cient\u00edfica, uso de ferramentas de busca, refer\u00eancias bibliogr\u00e1ficas e estruturas formais de divulga\u00e7\u00e3o cient\u00edfica. Desenvolvimento de um tema de pesquia ou projeto de engenharia, com o formato de um trabalho de inicia\u00e7\u00e3o cient\u00edfica, sob a orienta\u00e7\u00e3o de um professor ou pesquisador autorizado pela Comiss\u00e3o de Curso. Entrega e apresenta\u00e7\u00e3o de documento t\u00e9cnico no final da disciplina.\"\n\nReplace-Text \"Organization and formalism of the development of scientific work. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project at the end of the course.\" \"Organization and formalism of the development of scientific work or engineering project. Scientific writing techniques, use of search tools, bibliographic references and formal structures of scientific dissemination. Development of an individual research or engineering project topic, with the format of a Scientific Initiation work, under the guidance of a professor or researcher authorized by the Course Committee. Delivery and presentation of research project document at the end of the course.\"\n\nReplace-Text \"Aulas expositivas, reuni\u00f5es com professor orientador, desenvolvimento de projeto de pesquisa e elabora\u00e7\u00e3o de projeto de pesquisa.\" \"Aulas expositivas, reuni\u00f5es com professor orientador, desenvolvimento de projeto de pesquisa e/ou engenharia e elabora\u00e7\u00e3o de projeto de pesquisa e/ou engenharia. Visitas t\u00e9cnicas em institutos ou empresas da \u00e1rea cient\u00edfica ou de engenharia.\"\n"}
